$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 61.32777777777778
$ws.Range("E2").Value = 67.20555555555556
$ws.Range("G2").Value = -5.87777777777778
$ws.Range("D3").Value = 15.67975
$ws.Range("E3").Value = 18.26519444444445
$ws.Range("G3").Value = -2.585444444444446
$ws.Range("D4").Value = 153.7516388888889
$ws.Range("E4").Value = 162.8410555555556
$ws.Range("G4").Value = -9.089416666666665
$ws.Range("D5").Value = 6127.447361111112
$ws.Range("E5").Value = 7508.073027777777
$ws.Range("G5").Value = -1380.625666666665
$ws.Range("D9").Value = 102.475462962963
$ws.Range("E9").Value = 114.0638888888889
$ws.Range("G9").Value = -11.5884259259259
$ws.Range("D10").Value = 28.48948611111111
$ws.Range("E10").Value = 32.93819444444444
$ws.Range("G10").Value = -4.448708333333332
$ws.Range("D11").Value = 150.8274537037037
$ws.Range("E11").Value = 176.8305740740741
$ws.Range("G11").Value = -26.00312037037034
$ws.Range("D12").Value = 11873.96275462963
$ws.Range("E12").Value = 14621.79834259259
$ws.Range("G12").Value = -2747.835587962965
$ws.Range("D16").Value = 93.24490740740741
$ws.Range("E16").Value = 99.12037037037038
$ws.Range("G16").Value = -5.87546296296297
$ws.Range("D17").Value = 24.5097037037037
$ws.Range("E17").Value = 26.91130555555556
$ws.Range("G17").Value = -2.401601851851854
$ws.Range("D18").Value = 160.2035925925926
$ws.Range("E18").Value = 145.1042222222222
$ws.Range("G18").Value = 15.09937037037039
$ws.Range("D19").Value = 9635.487814814815
$ws.Range("E19").Value = 10822.44786111111
$ws.Range("G19").Value = -1186.960046296297
$ws.Range("D23").Value = 133.5893518518519
$ws.Range("E23").Value = 105.9212962962963
$ws.Range("G23").Value = 27.66805555555555
$ws.Range("D24").Value = 36.13469444444444
$ws.Range("E24").Value = 26.42299074074074
$ws.Range("G24").Value = 9.711703703703702
$ws.Range("D25").Value = 144.3176712962963
$ws.Range("E25").Value = 116.863537037037
$ws.Range("G25").Value = 27.45413425925926
$ws.Range("D26").Value = 13976.09596759259
$ws.Range("E26").Value = 9238.73998148148
$ws.Range("G26").Value = 4737.355986111112
$ws.Range("D30").Value = 234.2064814814815
$ws.Range("E30").Value = 99.17129629629629
$ws.Range("G30").Value = 135.0351851851852
$ws.Range("D31").Value = 59.10475
$ws.Range("E31").Value = 23.67728703703704
$ws.Range("G31").Value = 35.42746296296296
$ws.Range("D32").Value = 116.7894490740741
$ws.Range("E32").Value = 117.0621759259259
$ws.Range("F32").Value = "late > early"
$ws.Range("G32").Value = -0.2727268518518429
$ws.Range("D33").Value = 21756.63385648148
$ws.Range("E33").Value = 8837.139092592592
$ws.Range("G33").Value = 12919.49476388889
$ws.Range("D37").Value = 162.313888888889
$ws.Range("E37").Value = 151.0166666666667
$ws.Range("G37").Value = 11.29722222222227
$ws.Range("D38").Value = 42.22015277777778
$ws.Range("E38").Value = 30.12684259259259
$ws.Range("G38").Value = 12.0933101851852
$ws.Range("D39").Value = 129.5212314814815
$ws.Range("E39").Value = 60.95931481481481
$ws.Range("G39").Value = 68.56191666666666
$ws.Range("D40").Value = 16045.21368981482
$ws.Range("E40").Value = 9303.403685185185
$ws.Range("G40").Value = 6741.810004629631
